# Apply the "xlsx 2 json" convention-renaming edit:
#  - strip the "Obj-"/"Tbl-"/"Enum-" prefixes from sheet (tab) names
#  - rewrite each sheet's title cell (A1) to the "##OBJ--"/"##TBL--"/"##ENUM--" marker form
#  - swap the English field-name row and the Chinese label row so English names
#    come right after the title, and the type row carries "xid" as its key type
#  - add a new "regions" column (with its Chinese label + sample data) to the city sheet

$wb = $excel.ActiveWorkbook

# ---- 1. Sheet (tab) names --------------------------------------------------
$wsUser = $wb.Worksheets.Item(1)
$wsCity = $wb.Worksheets.Item(2)
$wsFood = $wb.Worksheets.Item(3)
$wsLang = $wb.Worksheets.Item(4)

$wsUser.Name = "user"
$wsCity.Name = "city"
$wsFood.Name = "food"
$wsLang.Name = "language"

# ---- 2. Sheet1 "user" -------------------------------------------------------
$wsUser.Range("A1").Value = "##OBJ--用户表"
$wsUser.Range("B2").Value = "type"

# ---- 3. Sheet2 "city" --------------------------------------------------------
$wsCity.Range("A1").Value = "##TBL--城市表"

# row2: English field names (+ new "regions" column)
$wsCity.Range("A2").Value = "id"
$wsCity.Range("B2").Value = "cityName"
$wsCity.Range("C2").Value = "province"
$wsCity.Range("D2").Value = "food"
$wsCity.Range("E2").Value = "regions"

# row3: type row now starts with "xid" (+ new "list:str" type for regions)
$wsCity.Range("A3").Value = "xid"
$wsCity.Range("E3").Value = "list:str"

# row4: Chinese labels (+ new "行政区" column)
$wsCity.Range("A4").Value = "索引"
$wsCity.Range("B4").Value = "城市名"
$wsCity.Range("C4").Value = "省份"
$wsCity.Range("D4").Value = "特产"
$wsCity.Range("E4").Value = "行政区"

# new data for the "regions" column
$wsCity.Range("E5").Value = "天河区,海珠区"
$wsCity.Range("E6").Value = "西湖区,新建区"

# ---- 4. Sheet3 "food" --------------------------------------------------------
$wsFood.Range("A1").Value = "##TBL--特产表"

# row2: English field names
$wsFood.Range("A2").Value = "id"
$wsFood.Range("B2").Value = "name"
$wsFood.Range("C2").Value = "type"
$wsFood.Range("D2").Value = "city"

# row3: type row now starts with "xid"
$wsFood.Range("A3").Value = "xid"

# row4: Chinese labels
$wsFood.Range("A4").Value = "索引"
$wsFood.Range("B4").Value = "食物名"
$wsFood.Range("C4").Value = "类别"
$wsFood.Range("D4").Value = "关联城市"

# ---- 5. Sheet4 "language" ----------------------------------------------------
$wsLang.Range("A1").Value = "##ENUM--语言表"
